# Update Atomos_Profits leve-crafting profit figures across all 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) following a scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 70019
$ws.Range("I21").Value = 70019
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 70019
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -69551
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 70019
$ws.Range("I23").Value = 70019
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 70019
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -69785
$ws.Range("N23").ClearContents()
$ws.Range("H33").Value = 432.94446
$ws.Range("I33").Value = 420.7143
$ws.Range("J33").Value = 475.75
$ws.Range("K33").Value = 420.7143
$ws.Range("L33").Value = 475.75
$ws.Range("M33").Value = -191.7143
$ws.Range("N33").Value = -933.75
$ws.Range("H63").Value = 39333.332
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 39333.332
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 39333.332
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -40581.332
$ws.Range("H66").Value = 39333.332
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 39333.332
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 117999.996
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -124239.996
$ws.Range("H76").Value = 3201
$ws.Range("I76").Value = 3089.8
$ws.Range("J76").Value = 3571.6667
$ws.Range("K76").Value = 3089.8
$ws.Range("L76").Value = 3571.6667
$ws.Range("M76").Value = -2774.8
$ws.Range("N76").Value = -4201.6667
$ws.Range("H79").Value = 3201
$ws.Range("I79").Value = 3089.8
$ws.Range("J79").Value = 3571.6667
$ws.Range("K79").Value = 3089.8
$ws.Range("L79").Value = 3571.6667
$ws.Range("M79").Value = -1997.8
$ws.Range("N79").Value = -5755.6667
$ws.Range("H87").Value = 25773.4
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 25773.4
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 25773.4
$ws.Range("N87").Value = -28269.4
$ws.Range("H90").Value = 25773.4
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 25773.4
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 77320.20000000001
$ws.Range("N90").Value = -89800.20000000001
$ws.Range("H113").Value = 3685.4583
$ws.Range("I113").Value = 3817.0833
$ws.Range("J113").Value = 3553.8333
$ws.Range("K113").Value = 3817.0833
$ws.Range("L113").Value = 3553.8333
$ws.Range("M113").Value = -563.0832999999998
$ws.Range("N113").Value = -10061.8333
$ws.Range("H137").Value = 4171094.8
$ws.Range("I137").Value = 7697944
$ws.Range("J137").Value = 3000.2727
$ws.Range("K137").Value = 23093832
$ws.Range("L137").Value = 9000.8181
$ws.Range("M137").Value = -23091282
$ws.Range("N137").Value = -14100.8181
$ws.Range("H138").Value = 3097.0532
$ws.Range("I138").Value = 1297.4857
$ws.Range("J138").Value = 4671.675
$ws.Range("K138").Value = 3892.4571
$ws.Range("L138").Value = 14015.025
$ws.Range("M138").Value = 1247.5429
$ws.Range("N138").Value = -24295.025

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7575.9
$ws.Range("I32").Value = 5197.5
$ws.Range("J32").Value = 20062.5
$ws.Range("K32").Value = 5197.5
$ws.Range("L32").Value = 20062.5
$ws.Range("M32").Value = -4910.5
$ws.Range("N32").Value = -20636.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17873
$ws.Range("I82").Value = 4371.4
$ws.Range("J82").Value = 34750
$ws.Range("K82").Value = 4371.4
$ws.Range("L82").Value = 34750
$ws.Range("M82").Value = -3988.4
$ws.Range("N82").Value = -35516
$ws.Range("H85").Value = 17873
$ws.Range("I85").Value = 4371.4
$ws.Range("J85").Value = 34750
$ws.Range("K85").Value = 4371.4
$ws.Range("L85").Value = 34750
$ws.Range("M85").Value = -3045.4
$ws.Range("N85").Value = -37402
$ws.Range("H99").Value = 1695.3889
$ws.Range("I99").Value = 1147.8
$ws.Range("J99").Value = 4433.3335
$ws.Range("K99").Value = 1147.8
$ws.Range("L99").Value = 4433.3335
$ws.Range("M99").Value = 350.2
$ws.Range("N99").Value = -7429.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 20001.25
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 20001.25
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 20001.25
$ws.Range("N13").Value = -20279.25
$ws.Range("H31").Value = 2783000.2
$ws.Range("I31").Value = 5557570
$ws.Range("J31").Value = 8430.5
$ws.Range("K31").Value = 5557570
$ws.Range("L31").Value = 8430.5
$ws.Range("M31").Value = -5557275
$ws.Range("N31").Value = -9020.5
$ws.Range("H34").Value = 2783000.2
$ws.Range("I34").Value = 5557570
$ws.Range("J34").Value = 8430.5
$ws.Range("K34").Value = 5557570
$ws.Range("L34").Value = 8430.5
$ws.Range("M34").Value = -5557368
$ws.Range("N34").Value = -8834.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1684.7941
$ws.Range("I131").Value = 1230
$ws.Range("J131").Value = 2043.8422
$ws.Range("K131").Value = 3690
$ws.Range("L131").Value = 6131.5266
$ws.Range("M131").Value = 1350
$ws.Range("N131").Value = -16211.5266

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2158.4707
$ws.Range("I102").Value = 1574.5
$ws.Range("J102").Value = 3560
$ws.Range("K102").Value = 1574.5
$ws.Range("L102").Value = 3560
$ws.Range("M102").Value = 47.5
$ws.Range("N102").Value = -6804
$ws.Range("H126").Value = 4184.5386
$ws.Range("I126").Value = 1999.6666
$ws.Range("J126").Value = 4840
$ws.Range("K126").Value = 5998.9998
$ws.Range("L126").Value = 14520
$ws.Range("M126").Value = -3528.9998
$ws.Range("N126").Value = -19460
$ws.Range("H132").Value = 3090.55
$ws.Range("I132").Value = 1899.8889
$ws.Range("J132").Value = 4064.7273
$ws.Range("K132").Value = 5699.6667
$ws.Range("L132").Value = 12194.1819
$ws.Range("M132").Value = -3169.6667
$ws.Range("N132").Value = -17254.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3925
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3925
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3925
$ws.Range("N40").Value = -4197
$ws.Range("H93").Value = 2403.9092
$ws.Range("I93").Value = 1529.875
$ws.Range("J93").Value = 4734.6665
$ws.Range("K93").Value = 1529.875
$ws.Range("L93").Value = 4734.6665
$ws.Range("M93").Value = -281.875
$ws.Range("N93").Value = -7230.6665
$ws.Range("H133").Value = 29500
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 29500
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 29500
$ws.Range("N133").Value = -34560
$ws.Range("H140").Value = 35000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 35000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 35000
$ws.Range("N140").Value = -45360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 70003.75
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 70003.75
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 70003.75
$ws.Range("N11").Value = -70287.75
$ws.Range("H122").Value = 455947.97
$ws.Range("I122").Value = 910060.4399999999
$ws.Range("J122").Value = 1835.4546
$ws.Range("K122").Value = 2730181.32
$ws.Range("L122").Value = 5506.3638
$ws.Range("M122").Value = -2727731.32
$ws.Range("N122").Value = -10406.3638
$ws.Range("H125").Value = 30000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 30000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840

